$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# "There is no need to implement the parser simply print a message stating
#  which parser would have parsed the file."
#   -> "There is no need to implement the parsing logic simply print a
#       message stating which parser would have parsed the file."
#
# Only the first occurrence of "parser" (immediately before "simply") is
# reworded to "parsing logic"; the later "which parser would have parsed"
# stays untouched.
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$anchor = "implement the parser simply"
$idx = $full.IndexOf($anchor)

if ($idx -ge 0) {
    # "parser" sits right after "implement the " inside the anchor text;
    # only its trailing "er" needs to become "ing logic" so that
    # "pars" + "ing logic" => "parsing logic".
    $parserStart = $idx + "implement the ".Length
    $erStart = $parserStart + "pars".Length
    $erEnd = $erStart + "er".Length

    $erRange = $d.Range($erStart, $erEnd)
    $erRange.Text = ""

    $insertPoint = $d.Range($erStart, $erStart)
    $insertPoint.InsertAfter("ing logic")

    # Word drops a "_GoBack" bookmark at the location of the last edit.
    $goBackPos = $erStart + "ing logic".Length
    $goBackRange = $d.Range($goBackPos, $goBackPos)
    try {
        $d.Bookmarks.Add("_GoBack", $goBackRange)
    } catch {
        # Non-fatal if this runtime's Bookmarks collection rejects the call;
        # the wording fix above is the substantive part of the edit.
    }
}
